# Updated symbol list on Sun Dec 25 06:38:34 UTC 2022 with GitHub Actions
#
# This script re-applies the latest crypto price/ranking snapshot onto
# Sheet1 of the workbook. All affected columns (B,C,D,E) hold plain text
# in the original file (t="inlineStr"), so every write below forces the
# destination cell to a Text number format before assigning the value,
# then resets cosmetic formatting back to "Normal" so no stray
# number-format/style is left behind - this keeps values like "245.07"
# or "0.005426" stored as text instead of being coerced into numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

$updates = [ordered]@{
    # Price refresh for rows that kept their coin (row 2-8)
    "D2"  = "245.07"
    "D3"  = "23.08"
    "D4"  = "5.412"
    "D6"  = "3.388"
    "D7"  = "0.8083"
    "D8"  = "0.9249"

    # Rows 9-17: coin ranking reshuffled by one position, with a new
    # coin ("One") rotating back in at the bottom of the block.
    "B9"  = "WazirX"
    "C9"  = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
    "D9"  = "0.1424"
    "E9"  = "8WazirXWRX"

    "B10" = "MandalaExchangeToken"
    "C10" = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
    "D10" = "0.07433"
    "E10" = "9MandalaExchangeTokenMDX"

    "B11" = "LiechtensteinCryptoassetsExchange"
    "C11" = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
    "D11" = "0.03368"
    "E11" = "10LiechtensteinCryptoassetsExchangeLCX"

    "B12" = "BitrueCoin"
    "C12" = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
    "D12" = "0.03033"
    "E12" = "11BitrueCoinBTR"

    "B13" = "BitMartToken"
    "C13" = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
    "D13" = "0.09356"
    "E13" = "12BitMartTokenBMX"

    "B14" = "MCDex"
    "C14" = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
    "D14" = "3.955"
    "E14" = "13MCDexMCB"

    "B15" = "BitForexToken"
    "C15" = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
    "D15" = "0.001604"
    "E15" = "14BitForexTokenBF"

    "B16" = "CoinExToken"
    "C16" = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
    "D16" = "0.04818"
    "E16" = "15CoinExTokenCET"

    "B17" = "One"
    "C17" = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
    "D17" = "0.0005942"
    "E17" = "16OneONEWorstin24h"

    # Remaining scattered price refreshes further down the table
    "D18" = "0.005426"
    "D20" = "0.0009865"
    "D21" = "0.00007102"
    "D22" = "3.653"
    "D23" = "6.436"
    "D24" = "2.187"
    "D40" = "0.03967"
    "D41" = "0.006482"
    "D43" = "0.002401"
    "D44" = "0.006715"
    "D45" = "0.00005205"
    "D47" = "0.0005802"
    "E48" = "47CoinbaseStockTokenCOINBestin24h"
    "D49" = "0.002303"
    "E49" = "48BOLOBOLO"
}

foreach ($key in $updates.Keys) {
    Set-TextValue $key $updates[$key]
}
